# Update the date line in the title paragraph.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2024-06-11 Tuesday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-12 Wednesday", 2)

# Update the division-problem answers in the table, addressed by
# (row, column) so the duplicate "963÷2=481, 1" cells are disambiguated.
$tbl = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; Old = "321÷2=160, 1"; New = "868÷5=173, 3" },
    @{ Row = 1;  Col = 2; Old = "536÷2=268, 0"; New = "694÷9=77, 1" },
    @{ Row = 1;  Col = 3; Old = "517÷6=86, 1";  New = "117÷6=19, 3" },
    @{ Row = 1;  Col = 4; Old = "275÷9=30, 5";  New = "646÷3=215, 1" },
    @{ Row = 1;  Col = 5; Old = "870÷3=290, 0"; New = "726÷3=242, 0" },

    @{ Row = 5;  Col = 1; Old = "637÷2=318, 1"; New = "392÷8=49, 0" },
    @{ Row = 5;  Col = 2; Old = "963÷2=481, 1"; New = "196÷7=28, 0" },
    @{ Row = 5;  Col = 3; Old = "119÷3=39, 2";  New = "941÷3=313, 2" },
    @{ Row = 5;  Col = 4; Old = "372÷8=46, 4";  New = "975÷9=108, 3" },
    @{ Row = 5;  Col = 5; Old = "699÷8=87, 3";  New = "565÷3=188, 1" },

    @{ Row = 9;  Col = 1; Old = "654÷6=109, 0"; New = "585÷6=97, 3" },
    @{ Row = 9;  Col = 2; Old = "242÷8=30, 2";  New = "952÷8=119, 0" },
    @{ Row = 9;  Col = 3; Old = "236÷6=39, 2";  New = "421÷7=60, 1" },
    @{ Row = 9;  Col = 4; Old = "397÷5=79, 2";  New = "534÷6=89, 0" },
    @{ Row = 9;  Col = 5; Old = "495÷2=247, 1"; New = "171÷8=21, 3" },

    @{ Row = 13; Col = 1; Old = "963÷2=481, 1"; New = "731÷3=243, 2" },
    @{ Row = 13; Col = 2; Old = "879÷7=125, 4"; New = "773÷8=96, 5" },
    @{ Row = 13; Col = 3; Old = "555÷8=69, 3";  New = "939÷7=134, 1" },
    @{ Row = 13; Col = 4; Old = "638÷4=159, 2"; New = "920÷2=460, 0" },
    @{ Row = 13; Col = 5; Old = "450÷6=75, 0";  New = "375÷9=41, 6" },

    @{ Row = 17; Col = 1; Old = "109÷8=13, 5"; New = "508÷2=254, 0" },
    @{ Row = 17; Col = 2; Old = "874÷6=145, 4"; New = "865÷8=108, 1" },
    @{ Row = 17; Col = 3; Old = "114÷5=22, 4";  New = "887÷4=221, 3" },
    @{ Row = 17; Col = 4; Old = "393÷9=43, 6";  New = "293÷7=41, 6" },
    @{ Row = 17; Col = 5; Old = "128÷9=14, 2";  New = "154÷7=22, 0" }
)

foreach ($u in $updates) {
    $cell = $tbl.Rows.Item($u.Row).Cells.Item($u.Col)
    $cell.Range.Text = $u.New
}

Write-Output "done"
